$d = $word.ActiveDocument
$wNs = ' xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) Insert four new paragraphs ("Skenario Pengguna" subsection) right before
#    the empty paragraph that follows "... Figma merupakan ….." and precedes
#    the "Product backlog" heading.
# ---------------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("Figma merupakan", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$figmaPara = $d.Range($anchor.Start, $anchor.Start).Paragraphs(1)
$targetPara = $figmaPara.Next()
$insertPoint = $d.Range($targetPara.Range.Start, $targetPara.Range.Start)

$rpr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$ppr = '<w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="360" w:firstLine="360"/><w:jc w:val="both"/>' + $rpr + '</w:pPr>'

$p1 = '<w:p' + $wNs + '>' + $ppr + '</w:p>'
$p2 = '<w:p' + $wNs + '>' + $ppr + '<w:r>' + $rpr + '<w:t>Skenario Pengguna</w:t></w:r></w:p>'
$p3 = '<w:p' + $wNs + '>' + $ppr + '</w:p>'
$testerText = 'Tester yang diperlukan pada proyek pengembangan ini ialah 9 orang dengan masing-masing user diwakili oleh 3 orang. Tester diperlukan untuk evaluasi lebih lanjut terkait rancangan proyek pengembangan yang telah dibuat. Feedback yang diberikan user akan dievaluasi lebih lanjut guna memperbaiki kesalahan yang dibuat pada pengembangan proyek terkait. '
$p4 = '<w:p' + $wNs + '>' + $ppr + '<w:r>' + $rpr + '<w:t xml:space="preserve">' + $testerText + '</w:t></w:r></w:p>'
# Trailing empty paragraph so InsertXML does not merge the "Tester..." text
# into the pre-existing (unrelated) paragraph that follows the insertion
# point; it is removed again right afterwards.
$p5 = '<w:p' + $wNs + '></w:p>'

$insertPoint.InsertXML($p1 + $p2 + $p3 + $p4 + $p5) | Out-Null

# Re-locate the freshly inserted paragraphs and drop the spacer paragraph
# that kept the real target paragraph from being merged/clobbered.
$anchor2 = $d.Content
$anchor2.Find.Execute("Figma merupakan", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$figmaPara2 = $d.Range($anchor2.Start, $anchor2.Start).Paragraphs(1)
$emptyPara1 = $figmaPara2.Next()
$skenarioPara = $emptyPara1.Next()
$emptyPara2 = $skenarioPara.Next()
$testerPara = $emptyPara2.Next()
$spacerPara = $testerPara.Next()
$spacerPara.Range.Delete() | Out-Null

# ---------------------------------------------------------------------------
# 2) Add a <w:lastRenderedPageBreak/> in front of the "Pada " tab-run that
#    introduces the product-backlog table.
# ---------------------------------------------------------------------------
$padaAnchor = $d.Content
$padaAnchor.Find.Execute("Pada proyek pengembangan perangkat lunak Manajemen Nilai Mahasiswa ini, tim Megantropus telah membuat product backlog seperti berikut", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$padaPara = $d.Range($padaAnchor.Start, $padaAnchor.Start).Paragraphs(1)
$padaRange = $d.Range($padaPara.Range.Start, $padaPara.Range.End)
$padaXml = '<w:p' + $wNs + ' w14:paraId="716EC69C" w:rsidR="006F1BE0" w:rsidRDefault="00E13A17" w:rsidP="00807DBA" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/>' + $rpr + '</w:pPr><w:r>' + $rpr + '<w:lastRenderedPageBreak/><w:tab/><w:t xml:space="preserve">Pada </w:t></w:r><w:r w:rsidR="00900549">' + $rpr + '<w:t>proyek pengembangan perangkat lunak Manajemen Nilai Mahasiswa ini, tim Megantropus telah membuat product backlog seperti berikut</w:t></w:r><w:r w:rsidR="006F1BE0">' + $rpr + '<w:t>.</w:t></w:r></w:p>'
$padaRange.InsertXML($padaXml) | Out-Null

# ---------------------------------------------------------------------------
# 3) Remove the (now superfluous) <w:lastRenderedPageBreak/> in the table
#    cell that just holds the literal "6" (row 7, column 1 of the backlog
#    table).
# ---------------------------------------------------------------------------
$table = $d.Tables.Item(1)
$cell = $table.Cell(7, 1)
$cellPara = $cell.Range.Paragraphs(1)
$cellRange = $d.Range($cellPara.Range.Start, $cellPara.Range.End)
$cellXml = '<w:p' + $wNs + ' w14:paraId="0F2BFFEF" w14:textId="0B185396" w:rsidR="00374CB6" w:rsidRDefault="00374CB6" w:rsidP="00807DBA" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/>' + $rpr + '</w:pPr><w:r>' + $rpr + '<w:t>6</w:t></w:r></w:p>'
$cellRange.InsertXML($cellXml) | Out-Null
